$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cell values (ifo qoq forecast-error recompute)
$ws.Range("J25").Value = -8.161966548722575
$ws.Range("K25").Value = -2.809361661215774
$ws.Range("I26").Value = -8.159011788180498
$ws.Range("J26").Value = -2.806510317809929
$ws.Range("H27").Value = -8.16335617673839
$ws.Range("I27").Value = -2.810749322399487
$ws.Range("G28").Value = -8.159011788180498
$ws.Range("H28").Value = -2.806510317809929
$ws.Range("F29").Value = -8.120522770801983
$ws.Range("G29").Value = -2.768016911160131
$ws.Range("H29").Value = -2.913197998145364
$ws.Range("I29").Value = -1.164222876508844
$ws.Range("J29").Value = 1.461604321912361
$ws.Range("K29").Value = -4.181974971976672
$ws.Range("E30").Value = -8.120540966007212
$ws.Range("F30").Value = -2.768035101504736
$ws.Range("G30").Value = -2.7982377905888
$ws.Range("H30").Value = -1.164222876508858
$ws.Range("I30").Value = 1.46160432191239
$ws.Range("J30").Value = -4.181974971976672
$ws.Range("D31").Value = -5.582366683264027
$ws.Range("E31").Value = -2.606510317809935
$ws.Range("F31").Value = -2.684726724817722
$ws.Range("G31").Value = -1.164222876508844
$ws.Range("H31").Value = 1.461604321912347
$ws.Range("I31").Value = -4.181974971976658
$ws.Range("C32").Value = -1.859011788180498
$ws.Range("D32").Value = 1.093489682190071
$ws.Range("E32").Value = -1.312938874122935
$ws.Range("F32").Value = -0.16422287650885
$ws.Range("G32").Value = 1.461604321912398
$ws.Range("H32").Value = -4.181974971976663
$ws.Range("B33").Value = -1.089459271323719
$ws.Range("C33").Value = 0.1211663704742572
$ws.Range("D33").Value = -1.478967116022048
$ws.Range("E33").Value = 0.1157683069622242
$ws.Range("F33").Value = 2.256329091196832
$ws.Range("G33").Value = -3.576252325814792
$ws.Range("H33").Value = 0.1742581798826135
$ws.Range("I33").Value = -0.5482747430672961
$ws.Range("J33").Value = -1.451141586996598
$ws.Range("K33").Value = 2.23808601891443
$ws.Range("B34").Value = -1.588722206811738
$ws.Range("C34").Value = -3.151932303024964
$ws.Range("D34").Value = 0.3674243044610539
$ws.Range("E34").Value = 4.3369229999808
$ws.Range("F34").Value = -2.715040054714407
$ws.Range("G34").Value = 0.4789198153989105
$ws.Range("H34").Value = -0.427397108966687
$ws.Range("I34").Value = -1.451141586996612
$ws.Range("J34").Value = 2.23808601891443
$ws.Range("B35").Value = -4.799062879723905
$ws.Range("C35").Value = 0.09539882002178501
$ws.Range("D35").Value = 5.294511079851064
$ws.Range("E35").Value = -2.259310354562728
$ws.Range("F35").Value = 0.4227867276758985
$ws.Range("G35").Value = -0.4314419592910641
$ws.Range("H35").Value = -1.306699781963914
$ws.Range("I35").Value = 2.238086018914458
$ws.Range("B36").Value = 0.812721528812574
$ws.Range("C36").Value = 4.953942165787652
$ws.Range("D36").Value = -2.956870662203044
$ws.Range("E36").Value = 0.3583149055470753
$ws.Range("F36").Value = -0.4477508806930921
$ws.Range("G36").Value = -1.315647347510321
$ws.Range("H36").Value = 2.238086018914458
$ws.Range("B37").Value = 0.9980800887900187
$ws.Range("C37").Value = -2.144503303453859
$ws.Range("D37").Value = 0.7785875935565794
$ws.Range("E37").Value = -0.2209348531114581
$ws.Range("F37").Value = -1.313788376955344
$ws.Range("G37").Value = 2.238086018914473
$ws.Range("H37").Value = 1.798459270573645
$ws.Range("I37").Value = 0.9080743123241604
$ws.Range("J37").Value = 2.174764871858827
$ws.Range("K37").Value = 2.581785915908512
$ws.Range("B38").Value = -4.427739562431512
$ws.Range("C38").Value = -0.5085356306463922
$ws.Range("D38").Value = 1.878549043769994
$ws.Range("E38").Value = 0.1640824532243532
$ws.Range("F38").Value = 3.052548475781066
$ws.Range("G38").Value = 1.993482555557406
$ws.Range("H38").Value = 0.9080743123241319
$ws.Range("I38").Value = 2.174764871858841
$ws.Range("J38").Value = 2.581785915908512
$ws.Range("B39").Value = -2.286871620153079
$ws.Range("C39").Value = -0.290770678289888
$ws.Range("D39").Value = -0.5572217042491578
$ws.Range("E39").Value = 2.249238898010802
$ws.Range("F39").Value = 1.782185565874443
$ws.Range("G39").Value = 0.9080743123241604
$ws.Range("H39").Value = 2.174764871858798
$ws.Range("I39").Value = 2.581785915908512
$ws.Range("B40").Value = 0.173268898576453
$ws.Range("C40").Value = -0.858249081784549
$ws.Range("D40").Value = 2.310430000587758
$ws.Range("E40").Value = 1.754521072957701
$ws.Range("F40").Value = 0.9080743123241319
$ws.Range("G40").Value = 2.174764871858827
$ws.Range("H40").Value = 2.581785915908512
$ws.Range("B41").Value = -2.153759632931141
$ws.Range("C41").Value = 1.96253099561207
$ws.Range("D41").Value = 0.8682245810140046
$ws.Range("E41").Value = 0.531039098243923
$ws.Range("F41").Value = 2.912659180991767
$ws.Range("G41").Value = 3.204872792212839
$ws.Range("H41").Value = 0.8261041740273922
$ws.Range("I41").Value = 0.9960589915708058
$ws.Range("J41").Value = 0.8671966487193004
$ws.Range("K41").Value = 0.474465624449067
$ws.Range("B42").Value = -0.6784352119512822
$ws.Range("C42").Value = 1.340740694025882
$ws.Range("D42").Value = 1.031288435998093
$ws.Range("E42").Value = 2.451533103558915
$ws.Range("F42").Value = 2.781785915908514
$ws.Range("G42").Value = 0.4082047373275373
$ws.Range("H42").Value = 0.8051410442068629
$ws.Range("I42").Value = 0.6764783301986341
$ws.Range("J42").Value = 0.2839473270420001
$ws.Range("B43").Value = -0.5250338534810766
$ws.Range("C43").Value = 1.123392218832621
$ws.Range("D43").Value = 2.760040968125835
$ws.Range("E43").Value = 3.013212894598098
$ws.Range("F43").Value = 0.6351697524705031
$ws.Range("G43").Value = 1.0737771510817
$ws.Range("H43").Value = 0.7914043845609768
$ws.Range("I43").Value = 0.2836599995723361
$ws.Range("B44").Value = 0.188257056176198
$ws.Range("C44").Value = 1.934899327679034
$ws.Range("D44").Value = 2.958115042475683
$ws.Range("E44").Value = 1.011685167471043
$ws.Range("F44").Value = 1.324838541200819
$ws.Range("G44").Value = 1.045575309115435
$ws.Range("H44").Value = 0.3513984122422613
$ws.Range("B45").Value = 1.944010268946087
$ws.Range("C45").Value = 2.331384361501406
$ws.Range("D45").Value = 0.08849833442987176
$ws.Range("E45").Value = 0.8372699328141238
$ws.Range("F45").Value = 0.7830637783740713
$ws.Range("G45").Value = 0.3530940214645995
$ws.Range("H45").Value = -0.1205515391266945
$ws.Range("I45").Value = -0.689832312574836
$ws.Range("B46").Value = 0.3493543950575599
$ws.Range("C46").Value = 0.3124941944969493
$ws.Range("D46").Value = 1.13466718197995
$ws.Range("E46").Value = 1.013226095514753
$ws.Range("F46").Value = 0.5835337946193699
$ws.Range("G46").Value = -0.01991609554966089
$ws.Range("H46").Value = -0.593701808683633
$ws.Range("B47").Value = -0.7067084814199234
$ws.Range("C47").Value = 0.475823087782075
$ws.Range("D47").Value = 0.9313544980887374
$ws.Range("E47").Value = 0.582365266164047
$ws.Range("F47").Value = -0.017682029428159
$ws.Range("G47").Value = -0.592866717998745
$ws.Range("B48").Value = 0.8219511122795922
$ws.Range("C48").Value = 1.06874533798873
$ws.Range("D48").Value = 0.750837591975525
$ws.Range("E48").Value = -0.2254573635392776
$ws.Range("F48").Value = -1.039202404151766
$ws.Range("B49").Value = -0.4155206925118478
$ws.Range("C49").Value = 0.04253363760471984
$ws.Range("D49").Value = -0.4474879771061924
$ws.Range("E49").Value = -0.9878138283037428
$ws.Range("B50").Value = -0.2385134798295212
$ws.Range("C50").Value = -0.6199905943205695
$ws.Range("D50").Value = -1.138587839195608
$ws.Range("B51").Value = -1.10660746326883
$ws.Range("C51").Value = -1.698025227524084
$ws.Range("B52").Value = -1.594906539899639

# Clear cells removed by the recompute (row now ends one column earlier)
$ws.Range("J45").ClearContents()
$ws.Range("I46").ClearContents()
$ws.Range("H47").ClearContents()
$ws.Range("G48").ClearContents()
$ws.Range("F49").ClearContents()
$ws.Range("E50").ClearContents()
$ws.Range("D51").ClearContents()
$ws.Range("C52").ClearContents()
$ws.Range("B53").ClearContents()
